$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-5 (columns D, J, K, L, M, N, O, P, Q)
# Row 2
$ws.Range("D2").Value = 44691
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 3000
$ws.Range("L2").Value = 3500
$ws.Range("M2").Value = 3250
$ws.Range("N2").Value = "$/docena de matas"
$ws.Range("O2").Value = "Región Metropolitana"
$ws.Range("P2").Value = 542
$ws.Range("Q2").Value = 6

# Row 3
$ws.Range("D3").Value = 44692
$ws.Range("J3").Value = 120

# Row 4
$ws.Range("D4").Value = 44687
$ws.Range("J4").Value = 160

# Row 5
$ws.Range("D5").Value = 44221
$ws.Range("J5").Value = 250
$ws.Range("K5").Value = 1300
$ws.Range("L5").Value = 1500
$ws.Range("M5").Value = 1420
$ws.Range("N5").Value = "$/atado"
$ws.Range("O5").Value = "Provincia de Diguillín"
$ws.Range("P5").Value = 1420
$ws.Range("Q5").Value = 1
